# Added games for 1/13/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Beat Vegas?" (column G) result for the 1/12/2021 games
#     that were already on the sheet (rows 54-59) ---
$ws.Range("G54").Value = "Yes"
$ws.Range("G55").Value = "Yes"
$ws.Range("G56").Value = "Yes"
$ws.Range("G57").Value = "No"
$ws.Range("G58").Value = "No"
$ws.Range("G59").Value = "Yes"

# --- Append the 1/13/2021 games (rows 60-67) ---
$newGames = @(
    @{ Row=60; Date=44209; Home="DET"; Away="MIL"; Spread=10.5;  Pred=23.3;  Diff=-12.8 },
    @{ Row=61; Date=44209; Home="CHO"; Away="DAL"; Spread=4;     Pred=-0.2;  Diff=4.2 },
    @{ Row=62; Date=44209; Home="NYK"; Away="BRK"; Spread=6;     Pred=26.5;  Diff=-20.5 },
    @{ Row=63; Date=44209; Home="MIN"; Away="MEM"; Spread=-3.5;  Pred=-2.3;  Diff=-1.2 },
    @{ Row=64; Date=44209; Home="OKC"; Away="LAL"; Spread=9;     Pred=8.1;   Diff=0.9 },
    @{ Row=65; Date=44209; Home="PHO"; Away="ATL"; Spread=-5.5;  Pred=-7.6;  Diff=2.1 },
    @{ Row=66; Date=44209; Home="LAC"; Away="NOP"; Spread=-6;    Pred=-4.8;  Diff=-1.2 },
    @{ Row=67; Date=44209; Home="SAC"; Away="POR"; Spread=4.5;   Pred=-3;    Diff=7.5 }
)

foreach ($g in $newGames) {
    $r = $g.Row
    $ws.Range("A$r").Value = $g.Date
    $ws.Range("A$r").NumberFormat = "yyyy\-mm\-dd"
    $ws.Range("B$r").Value = $g.Home
    $ws.Range("C$r").Value = $g.Away
    $ws.Range("D$r").Value = $g.Spread
    $ws.Range("E$r").Value = $g.Pred
    $ws.Range("F$r").Value = $g.Diff
}

# --- Restore the view state (scroll position / active selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$ws.Range("G71").Select() | Out-Null
